$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 909, shifting the existing rows 909:950 down to 910:951.
$ws.Rows.Item(909).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds the date as literal text (matching the rest of the
# column), so force text formatting before assigning the value to stop
# Excel from auto-converting the string into a date serial number, then
# restore the default "Normal" style so no stray number format sticks
# to the cell (the other cells in the column carry no explicit style).
$dateCell = $ws.Cells.Item(909, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/03/01"
$dateCell.Style = "Normal"

$ws.Cells.Item(909, 2).Value = "日"
$ws.Cells.Item(909, 3).Value = 13
$ws.Cells.Item(909, 4).Value = 201
